{"js": "// Renumber the PERSON_* placeholders that shifted down by one (firstname\n// inference fixes for contracts 18-24). Each old label is a unique,\n// complete run of text, so a simple search+replace per pair is safe and\n// order independent (no replacement's output text collides with another\n// pair's search text).\nconst replacements = [\n  [\"MUDr. [[PERSON_8]],\", \"MUDr. [[PERSON_4]],\"],\n  [\"sestry [[PERSON_9]],\", \"sestry [[PERSON_8]],\"],\n  [\"laborantky [[PERSON_10]].\", \"laborantky [[PERSON_9]].\"],\n  [\"asistentkou [[PERSON_11]],\", \"asistentkou [[PERSON_10]],\"],\n  [\"sonografistkou Mgr. [[PERSON_12]].\", \"sonografistkou Mgr. [[PERSON_11]].\"],\n  [\"MUDr. [[PERSON_13]],\", \"MUDr. [[PERSON_12]],\"],\n  [\"sestra [[PERSON_14]].\", \"sestra [[PERSON_13]].\"],\n  [\"provedla MUDr. [[PERSON_15]], radiolo\u017eka\", \"provedla MUDr. [[PERSON_14]], radiolo\u017eka\"],\n  [\"asistent: Bc. [[PERSON_16]]\", \"asistent: Bc. [[PERSON_15]]\"],\n  [\"l\u00e9ka\u0159: MUDr. [[PERSON_17]]\", \"l\u00e9ka\u0159: MUDr. [[PERSON_16]]\"],\n  [\"l\u00e9ka\u0159ka: MUDr. [[PERSON_18]]\", \"l\u00e9ka\u0159ka: MUDr. [[PERSON_17]]\"],\n  [\"technik: [[PERSON_19]]\", \"technik: [[PERSON_18]]\"],\n  [\"Mgr. [[PERSON_20]],\", \"Mgr. [[PERSON_19]],\"],\n  [\"Bc. [[PERSON_21]],\", \"Bc. [[PERSON_20]],\"],\n  [\"Mgr. [[PERSON_22]].\", \"Mgr. [[PERSON_21]].\"],\n  [\"cvi\u010den\u00ed dle metody DNS (pod dohledem Mgr. [[PERSON_23]]),\", \"cvi\u010den\u00ed dle metody DNS (pod dohledem Mgr. [[PERSON_22]]),\"],\n  [\"MUDr. [[PERSON_24]],\", \"MUDr. [[PERSON_23]],\"],\n  [\"sestra [[PERSON_25]].\", \"sestra [[PERSON_24]].\"],\n  [\"Interna A \u2014 prim\u00e1\u0159 MUDr. [[PERSON_26]],\", \"Interna A \u2014 prim\u00e1\u0159 MUDr. [[PERSON_25]],\"],\n  [\"Gynekologie \u2014 prim\u00e1\u0159ka MUDr. [[PERSON_27]],\", \"Gynekologie \u2014 prim\u00e1\u0159ka MUDr. [[PERSON_26]],\"],\n  [\"ORL \u2014 prim\u00e1\u0159 MUDr. [[PERSON_28]],\", \"ORL \u2014 prim\u00e1\u0159 MUDr. [[PERSON_27]],\"],\n  [\"Dermatologie \u2014 garant MUDr. [[PERSON_29]].\", \"Dermatologie \u2014 garant MUDr. [[PERSON_28]].\"],\n  [\"MUDr. [[PERSON_30]],\", \"MUDr. [[PERSON_29]],\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Renumber the PERSON_* placeholders that shifted down by one (firstname\n# inference fixes for contracts 18-24). Each old label below is a unique,\n# complete run of text, so a plain Find/Replace per pair is safe and\n# order independent (no replacement's output text collides with another\n# pair's search text).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"MUDr. [[PERSON_8]],\", \"MUDr. [[PERSON_4]],\"),\n    @(\"sestry [[PERSON_9]],\", \"sestry [[PERSON_8]],\"),\n    @(\"laborantky [[PERSON_10]].\", \"laborantky [[PERSON_9]].\"),\n    @(\"asistentkou [[PERSON_11]],\", \"asistentkou [[PERSON_10]],\"),\n    @(\"sonografistkou Mgr. [[PERSON_12]].\", \"sonografistkou Mgr. [[PERSON_11]].\"),\n    @(\"MUDr. [[PERSON_13]],\", \"MUDr. [[PERSON_12]],\"),\n    @(\"sestra [[PERSON_14]].\", \"sestra [[PERSON_13]].\"),\n    @(\"provedla MUDr. [[PERSON_15]], radiolo\u017eka\", \"provedla MUDr. [[PERSON_14]], radiolo\u017eka\"),\n    @(\"asistent: Bc. [[PERSON_16]]\", \"asistent: Bc. [[PERSON_15]]\"),\n    @(\"l\u00e9ka\u0159: MUDr. [[PERSON_17]]\", \"l\u00e9ka\u0159: MUDr. [[PERSON_16]]\"),\n    @(\"l\u00e9ka\u0159ka: MUDr. [[PERSON_18]]\", \"l\u00e9ka\u0159ka: MUDr. [[PERSON_17]]\"),\n    @(\"technik: [[PERSON_19]]\", \"technik: [[PERSON_18]]\"),\n    @(\"Mgr. [[PERSON_20]],\", \"Mgr. [[PERSON_19]],\"),\n    @(\"Bc. [[PERSON_21]],\", \"Bc. [[PERSON_20]],\"),\n    @(\"Mgr. [[PERSON_22]].\", \"Mgr. [[PERSON_21]].\"),\n    @(\"cvi\u010den\u00ed dle metody DNS (pod dohledem Mgr. [[PERSON_23]]),\", \"cvi\u010den\u00ed dle metody DNS (pod dohledem Mgr. [[PERSON_22]]),\"),\n    @(\"MUDr. [[PERSON_24]],\", \"MUDr. [[PERSON_23]],\"),\n    @(\"sestra [[PERSON_25]].\", \"sestra [[PERSON_24]].\"),\n    @(\"Interna A \u2014 prim\u00e1\u0159 MUDr. [[PERSON_26]],\", \"Interna A \u2014 prim\u00e1\u0159 MUDr. [[PERSON_25]],\"),\n    @(\"Gynekologie \u2014 prim\u00e1\u0159ka MUDr. [[PERSON_27]],\", \"Gynekologie \u2014 prim\u00e1\u0159ka MUDr. [[PERSON_26]],\"),\n    @(\"ORL \u2014 prim\u00e1\u0159 MUDr. [[PERSON_28]],\", \"ORL \u2014 prim\u00e1\u0159 MUDr. [[PERSON_27]],\"),\n    @(\"Dermatologie \u2014 garant MUDr. [[PERSON_29]].\", \"Dermatologie \u2014 garant MUDr. [[PERSON_28]].\"),\n    @(\"MUDr. [[PERSON_30]],\", \"MUDr. [[PERSON_29]],\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 0=wdFindContinue, 2=wdReplaceAll\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
